$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.019.39'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").Value = '2.906.38'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '364.66'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '103.09'
$ws.Range("E6").Value = '  -6.51%  '
$ws.Range("E7").Value = '  -5.21%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  -7.17%  '
$ws.Range("D10").Value = '36.89'
$ws.Range("E10").Value = '  -5.68%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("E12").Value = '  -4.51%  '
$ws.Range("D13").Value = '18.38'
$ws.Range("E13").Value = '  -6.12%  '
$ws.Range("D14").Value = '3.367.90'
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = '7.33'
$ws.Range("E15").Value = '  -5.74%  '
$ws.Range("D16").Value = '2.902.09'
$ws.Range("E16").Value = '  -1.55%  '
$ws.Range("D17").Value = '0.949'
$ws.Range("E17").Value = '  -4.00%  '
$ws.Range("D18").Value = '51.003.12'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.30'
$ws.Range("E19").Value = '  -6.66%  '
$ws.Range("D20").Value = '7.23'
$ws.Range("E20").Value = '  -4.70%  '
$ws.Range("D21").Value = '12.98'
$ws.Range("E21").Value = '  -6.91%  '
$ws.Range("D22").Value = '0.0₃0944'
$ws.Range("E22").Value = '  -4.19%  '
$ws.Range("D23").Value = '68.04'
$ws.Range("E23").Value = '  -3.64%  '
$ws.Range("D24").Value = '259.71'
$ws.Range("E24").Value = '  -3.41%  '
$ws.Range("D25").Value = '2.69'
$ws.Range("E25").Value = '  -4.60%  '
$ws.Range("D26").Value = '4.34'
$ws.Range("E26").Value = '  +3.75%  '
$ws.Range("E27").Value = '  -5.45%  '
$ws.Range("D29").Value = '25.86'
$ws.Range("E29").Value = '  -4.46%  '
$ws.Range("D30").Value = '7.31'
$ws.Range("E30").Value = '  -5.85%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").Value = '9.92'
$ws.Range("E32").Value = '  -5.74%  '
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("E34").Value = '  -6.47%  '
$ws.Range("D35").Value = '35.02'
$ws.Range("E35").Value = '  -6.92%  '
$ws.Range("D36").Value = '50.51'
$ws.Range("E36").Value = '  -3.50%  '
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '0.0421'
$ws.Range("E38").Value = '  -5.21%  '
$ws.Range("D39").Value = '2.79'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").Value = '3.12'
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("D41").Value = '16.85'
$ws.Range("E41").Value = '  -7.97%  '
$ws.Range("E42").Value = '  -7.30%  '
$ws.Range("D43").Value = '0.113'
$ws.Range("E43").Value = '  -5.48%  '
$ws.Range("D44").Value = '22.32'
$ws.Range("E44").Value = '  -3.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '117.60'
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").Value = '2.061.43'
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").Value = '3.19'
$ws.Range("E48").Value = '  -8.02%  '
$ws.Range("E49").Value = '  -8.34%  '
$ws.Range("D50").Value = '3.202.37'
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("E51").Value = '  -6.96%  '
